$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:BBETF0432"
$ws.Range("C2").Value = "NSE:AARTISURF"
$ws.Range("D2").Value = "NSE:SAIL"
$ws.Range("E2").Value = "NSE:DIXON"
$ws.Range("F2").Value = "NSE:HUDCO"

# Row 3
$ws.Range("C3").Value = "NSE:AIRAN"
$ws.Range("E3").Value = "NSE:HAL"

# Row 4
$ws.Range("B4").Value = "NSE:GRAPHITE"
$ws.Range("C4").Value = "NSE:AWHCL"
$ws.Range("E4").Value = "NSE:PATANJALI"

# Row 5
$ws.Range("B5").Value = "NSE:HINDCOPPER"
$ws.Range("C5").Value = "NSE:CROMPTON"
$ws.Range("E5").Value = "NSE:PAYTM"

# Row 6
$ws.Range("B6").Value = "NSE:HUDCO"
$ws.Range("C6").Value = "NSE:DCXINDIA"
$ws.Range("E6").Value = "NSE:PERSISTENT"

# Row 7 (E7 cleared to blank text)
$ws.Range("B7").Value = "NSE:JAIBALAJI"
$ws.Range("C7").Value = "NSE:HERANBA"
$ws.Range("E7").Value = "'"
$ws.Range("E7").Style = "Normal"

# Row 8 (E8 cleared to blank text)
$ws.Range("B8").Value = "NSE:LICNETFGSC"
$ws.Range("C8").Value = "NSE:IOLCP"
$ws.Range("E8").Value = "'"
$ws.Range("E8").Style = "Normal"

# Row 9 (E9 cleared to blank text)
$ws.Range("B9").Value = "NSE:MAITHANALL"
$ws.Range("C9").Value = "NSE:ITETF"
$ws.Range("E9").Value = "'"
$ws.Range("E9").Style = "Normal"

# Row 10 (E10 cleared to blank text)
$ws.Range("B10").Value = "NSE:MMTC"
$ws.Range("C10").Value = "NSE:JTLIND"
$ws.Range("E10").Value = "'"
$ws.Range("E10").Style = "Normal"

# Row 11 (E11 cleared to blank text)
$ws.Range("B11").Value = "NSE:MONQ50"
$ws.Range("C11").Value = "NSE:KAUSHALYA"
$ws.Range("E11").Value = "'"
$ws.Range("E11").Style = "Normal"

# Row 12 (E12 cleared to blank text)
$ws.Range("B12").Value = "NSE:PRAKASHSTL"
$ws.Range("C12").Value = "NSE:KOKUYOCMLN"
$ws.Range("E12").Value = "'"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "NSE:PTCIL"
$ws.Range("C13").Value = "NSE:NAVA"

# Row 14 (B14 cleared to blank text)
$ws.Range("B14").Value = "'"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "NSE:PONNIERODE"

# Row 15 (B15 cleared to blank text)
$ws.Range("B15").Value = "'"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "NSE:PRSMJOHNSN"

# Row 16 (B16 cleared to blank text)
$ws.Range("B16").Value = "'"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "NSE:PURVA"

# Row 17
$ws.Range("C17").Value = "NSE:SAKAR"

# Remove rows 18-29 (shrinks dimension to A1:F17)
$ws.Range("A18:F29").EntireRow.Delete()
